# Adds the new Lab Instructor / Lab Assistant staff rows (15-18) to the
# "007" worksheet, matching the "Some Endpoint, Insertion and json are Added"
# commit: four new people appended below the existing faculty list with
# Name / Designation / Photo / unique_id only.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the formatting already applied to the "Photo" column of an existing
# data row as the template for the new Photo cells (I15:I18).
$photoTemplate = $ws.Cells.Item(3, 9)

# Row data: Name, Designation, Photo URL (note trailing newline, matching
# the existing photo-url strings), unique_id
$rows = @(
    @{ Row = 15; Name = "GNANA SEKAR";   Designation = "Lab Instructor"; Photo = "/static/images/profile_photos/007/VEC-007-05-8.webp`n";  Id = "VEC-007-05-8" },
    @{ Row = 16; Name = "RAJENDRAN.P";   Designation = "Lab Instructor"; Photo = "/static/images/profile_photos/007/VEC-007-05-9.webp`n";  Id = "VEC-007-05-9" },
    @{ Row = 17; Name = "RANGARAJAN";    Designation = "Lab Instructor"; Photo = "/static/images/profile_photos/007/VEC-007-05-10.webp`n"; Id = "VEC-007-05-10" },
    @{ Row = 18; Name = "SELVAKUMAR.P";  Designation = "Lab Assistant";  Photo = "/static/images/profile_photos/007/VEC-007-05-19.webp`n"; Id = "VEC-007-05-19" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $nameCell = $ws.Cells.Item($rowNum, 1)
    $nameCell.Value = $r.Name
    $nameCell.Style = "Normal"
    $nameCell.Borders.LineStyle = 1

    $desigCell = $ws.Cells.Item($rowNum, 2)
    $desigCell.Value = $r.Designation
    $desigCell.Style = "Normal"
    $desigCell.Borders.LineStyle = 1
    $desigCell.Font.Color = 0

    $photoCell = $ws.Cells.Item($rowNum, 9)
    $photoCell.Value = $r.Photo
    $photoTemplate.Copy() | Out-Null
    $photoCell.PasteSpecial(-4122) | Out-Null

    $idCell = $ws.Cells.Item($rowNum, 10)
    $idCell.Value = $r.Id
    $idCell.Style = "Normal"
}

$excel.CutCopyMode = $false

[void]$ws.Range("J16").Select()
